# Applies the cryptocurrency price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that LOOKS like a plain number (e.g. '0.830', '9.00')
# while keeping it stored as literal text, matching the source data, which
# otherwise would be auto-converted to a number and lose formatting such as
# trailing zeros (Excel's normal text-to-number inference on cell assignment).
function Set-TextCell($ref, $value) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = '@'
    $cell.Value = $value
    $cell.ClearFormats()
}

$ws.Range('D2').Value = '94.273.57'
$ws.Range('E2').Value = '  +2.93%  '
$ws.Range('D3').Value = '3.103.69'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  -0.05%  '
Set-TextCell 'D5' '237.91'
$ws.Range('E5').Value = '  -2.02%  '
Set-TextCell 'D6' '616.59'
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('E7').Value = '  +2.58%  '
$ws.Range('E8').Value = '  -1.31%  '
$ws.Range('E9').Value = '  -0.06%  '
Set-TextCell 'D10' '0.830'
$ws.Range('E10').Value = '  +13.06%  '
$ws.Range('D11').Value = '3.100.80'
$ws.Range('E11').Value = '  -0.21%  '
Set-TextCell 'D12' '0.197'
$ws.Range('E12').Value = '  -2.45%  '
Set-TextCell 'D13' '0.0000244'
$ws.Range('E13').Value = '  -3.28%  '
$ws.Range('D14').Value = '93.755.35'
$ws.Range('E14').Value = '  +1.86%  '
Set-TextCell 'D15' '34.94'
$ws.Range('E15').Value = '  +1.18%  '
Set-TextCell 'D16' '5.41'
$ws.Range('E16').Value = '  -2.05%  '
$ws.Range('D17').Value = '3.680.62'
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('D18').Value = '3.098.74'
$ws.Range('E18').Value = '  -1.88%  '
Set-TextCell 'D19' '3.65'
$ws.Range('E19').Value = '  +0.39%  '
Set-TextCell 'D20' '14.87'
$ws.Range('E20').Value = '  +0.50%  '
Set-TextCell 'D21' '5.97'
$ws.Range('E21').Value = '  +2.67%  '
Set-TextCell 'D22' '445.65'
$ws.Range('E22').Value = '  -0.36%  '
Set-TextCell 'D23' '0.0000200'
$ws.Range('E23').Value = '  -1.13%  '
Set-TextCell 'D24' '9.00'
$ws.Range('E24').Value = '  -3.22%  '
Set-TextCell 'D25' '8.18'
$ws.Range('E25').Value = '  +3.71%  '
$ws.Range('E26').Value = '  -0.16%  '
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 'D27' '12.20'
$ws.Range('E27').Value = '  +4.67%  '
$ws.Range('B28').Value = 'Litecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 'D28' '86.07'
$ws.Range('E28').Value = '  -1.08%  '
$ws.Range('D29').Value = '3.269.45'
$ws.Range('E29').Value = '  -0.28%  '
$ws.Range('E30').Value = '  +0.22%  '
Set-TextCell 'D31' '0.244'
$ws.Range('E31').Value = '  +6.87%  '
Set-TextCell 'D32' '0.180'
$ws.Range('E32').Value = '  +7.74%  '
$ws.Range('E33').Value = '  -10.36%  '
Set-TextCell 'D34' '9.27'
$ws.Range('E34').Value = '  -0.54%  '
$ws.Range('E35').Value = '  +0.08%  '
Set-TextCell 'D36' '0.165'
$ws.Range('E36').Value = '  -2.78%  '
Set-TextCell 'D37' '7.90'
$ws.Range('E37').Value = '  -1.23%  '
Set-TextCell 'D38' '26.05'
$ws.Range('E38').Value = '  -0.87%  '
$ws.Range('E39').Value = '  -1.49%  '
Set-TextCell 'D40' '0.452'
$ws.Range('E40').Value = '  +4.31%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell 'D41' '477.07'
$ws.Range('E41').Value = '  -0.86%  '
Set-TextCell 'D42' '23.97'
$ws.Range('E42').Value = '  +8.12%  '
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell 'D43' '1.28'
$ws.Range('E43').Value = '  -2.06%  '
$ws.Range('B44').Value = 'MantraDAO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-TextCell 'D44' '3.76'
$ws.Range('E44').Value = '  -8.30%  '
Set-TextCell 'D45' '3.24'
$ws.Range('E45').Value = '  -6.51%  '
Set-TextCell 'D47' '160.88'
$ws.Range('E47').Value = '  +1.30%  '
Set-TextCell 'D48' '0.688'
$ws.Range('E48').Value = '  -1.66%  '
$ws.Range('E49').Value = '  -2.84%  '
Set-TextCell 'D50' '4.43'
$ws.Range('E50').Value = '  +1.31%  '
$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell 'D51' '43.80'
$ws.Range('E51').Value = '  -0.26%  '
